$d = $word.ActiveDocument

$d.Paragraphs(1).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t>F</w:t></w:r><w:r><w:t>ag Handling Ideas</w:t></w:r></w:p>')
$d.Paragraphs(2).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>F</w:t></w:r><w:r><w:t>ag Interpretation Using Reflection</w:t></w:r></w:p>')
$d.Paragraphs(3).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">The BrowserMonkey program needs to interpret a large number of tags and do this process in a modular way </w:t></w:r><w:r><w:t>to allow for future usage of the tag handling system.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Most</w:t></w:r><w:r><w:t xml:space="preserve"> importantly</w:t></w:r><w:r><w:t xml:space="preserve"> this must be done very</w:t></w:r><w:r><w:t xml:space="preserve"> efficiently. </w:t></w:r><w:r><w:t>It would be possible</w:t></w:r><w:r><w:t xml:space="preserve"> create a html </w:t></w:r><w:r><w:t>tag handler</w:t></w:r><w:r><w:t xml:space="preserve"> by using a</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>enormous</w:t></w:r><w:r><w:t xml:space="preserve"> if-else statement but this would be hard to debug and horrible to reuse or update (for example to new HTML standards). This is why we have decided to use reflection.</w:t></w:r></w:p>')
$d.Paragraphs(5).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">The program </w:t></w:r><w:r><w:t>can</w:t></w:r><w:r><w:t xml:space="preserve"> use</w:t></w:r><w:r><w:t xml:space="preserve"> an external file (such as a config file) that is easily editable to lookup required information for use while the program is running. An</w:t></w:r><w:r><w:t xml:space="preserve"> alternative</w:t></w:r><w:r><w:t xml:space="preserve"> way of doing this is having the information in the file loaded into the program when is it initially run.</w:t></w:r><w:r><w:t xml:space="preserve"> This allows for the addition of useful modular features to a program.</w:t></w:r></w:p>')
$d.Paragraphs(6).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t>Reflection</w:t></w:r><w:r><w:t xml:space="preserve"> is the process by which a computer program can observe and modify its own structure and behaviour.</w:t></w:r></w:p>')
$d.Paragraphs(10).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">Now I will give a brief introduction of how these techniques </w:t></w:r><w:r><w:t>could</w:t></w:r><w:r><w:t xml:space="preserve"> be implemented in Java.</w:t></w:r></w:p>')
$d.Paragraphs(28).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:pStyle w:val="Code"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">        }</w:t></w:r></w:p>')
$d.Paragraphs(30).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:color w:val="4F81BD" w:themeColor="accent1"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">This Code will allow you to read each line in a file and </w:t></w:r><w:r><w:t>execute the required code on each line</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>')
$d.Paragraphs(51).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t>This</w:t></w:r><w:r><w:t xml:space="preserve"> code reads the command line arguments and compares the first argument to all available classes in java if it finds a matching class </w:t></w:r><w:r><w:t>it will output a toString of each method that class contains.</w:t></w:r><w:r><w:t xml:space="preserve"> If there is no matching class it will throw an error.</w:t></w:r></w:p>')
$d.Paragraphs(53).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t xml:space="preserve">To apply the above methods in the BrowserMonkey Browser we </w:t></w:r><w:r><w:t>will</w:t></w:r><w:r><w:t xml:space="preserve"> use a Class for each HTML tag that contains code that can be used for painting the </w:t></w:r><w:r><w:t xml:space="preserve">component related to the current tag. </w:t></w:r><w:r><w:t xml:space="preserve">One option is to implement the system as an </w:t></w:r><w:r><w:t xml:space="preserve">abstract class or </w:t></w:r><w:r><w:t xml:space="preserve">use </w:t></w:r><w:r><w:t xml:space="preserve">an interface </w:t></w:r><w:r><w:t>to increase the modularity of the implementation.</w:t></w:r></w:p>')
$d.Paragraphs(54).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t>There will be an external file that we can use to relate the tag we’re looking at to the</w:t></w:r><w:r><w:t xml:space="preserve"> name of the equivalent mini Tag Class</w:t></w:r><w:r><w:t xml:space="preserve"> then we will be using reflection to get the necessary paint methods over to the renderer so that it can build the required component based on the tag.</w:t></w:r></w:p>')
$d.Paragraphs(56).Range.InsertXML('<w:p xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:r><w:t>A simple pseudo code representation of how the above ideas would work within the program</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>')
